# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" right after "总计" (before "2022-Q2"),
#   seeded from the "2022-Q2" sheet's layout/styling, then overwritten with
#   the new quarter's fund data (4 rows instead of 2).
# - Update the "总计" summary sheet: new first data row for 2022-Q4, and the
#   previously-existing rows shift down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the "2022-Q4" sheet by copying "2022-Q2" (keeps header/style),
#    then move it so the tab order is 总计, 2022-Q4, 2022-Q2, 2021-Q4.
#    NOTE: worksheet handles in this host are positional, not stable
#    identities -- any handle captured before a Move()/reorder can silently
#    resolve to a *different* sheet afterwards. Always re-fetch by Name
#    right after reordering, and do all the data writes after that.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $template)
$wb.Worksheets.Item("2022-Q2 (2)").Name = "2022-Q4"

$anchor = $wb.Worksheets.Item("2022-Q2")
$wb.Worksheets.Item("2022-Q4").Move($anchor, $null)

$newSheet = $wb.Worksheets.Item("2022-Q4")

# Helper: force text type for numeric-looking strings by flipping to a text
# number format, assigning the value, then clearing the format again so the
# cell ends up back at the default style (matches how the sibling data
# cells are stored: text with no explicit style) while keeping the
# "string" cell type instead of being silently coerced to a number.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Grow from 2 data rows to 4 data rows: duplicate row 3's formatting twice
# (so the new rows 4 & 5 carry the same index-column style as rows 2 & 3).
$newSheet.Range("A3:H3").Copy($newSheet.Range("A4:H4"))
$newSheet.Range("A3:H3").Copy($newSheet.Range("A5:H5"))

# Row 2
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "013166"
Set-TextValue $newSheet.Range("C2") "东兴宸祥量化混合A"
Set-TextValue $newSheet.Range("D2") "0.38"
Set-TextValue $newSheet.Range("E2") "93.88"
Set-TextValue $newSheet.Range("F2") "1.08"
Set-TextValue $newSheet.Range("G2") "0.0041"
$newSheet.Range("H2").Value = 7

# Row 3
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "009327"
Set-TextValue $newSheet.Range("C3") "东兴兴晟混合A"
Set-TextValue $newSheet.Range("D3") "0.38"
Set-TextValue $newSheet.Range("E3") "79.79"
Set-TextValue $newSheet.Range("F3") "0.98"
Set-TextValue $newSheet.Range("G3") "0.0037"
$newSheet.Range("H3").Value = 8

# Row 4
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "013167"
Set-TextValue $newSheet.Range("C4") "东兴宸祥量化混合C"
Set-TextValue $newSheet.Range("D4") "0.08"
Set-TextValue $newSheet.Range("E4") "93.88"
Set-TextValue $newSheet.Range("F4") "1.08"
Set-TextValue $newSheet.Range("G4") "0.0009"
$newSheet.Range("H4").Value = 7

# Row 5
$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet.Range("B5") "009328"
Set-TextValue $newSheet.Range("C5") "东兴兴晟混合C"
Set-TextValue $newSheet.Range("D5") "0.08"
Set-TextValue $newSheet.Range("E5") "79.79"
Set-TextValue $newSheet.Range("F5") "0.98"
Set-TextValue $newSheet.Range("G5") "0.0008"
$newSheet.Range("H5").Value = 8

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: push existing rows down by one and
#    insert the new 2022-Q4 totals as the new row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

$summary.Range("A2").Value = 0
Set-TextValue $summary.Range("B2") "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.01

$summary.Range("A3").Value = 1
Set-TextValue $summary.Range("B3") "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0

$summary.Range("A4").Value = 2
Set-TextValue $summary.Range("B4") "2021-Q4"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.05

# Restore the originally-active sheet/selection (the copy/move above left
# "2022-Q2" focused).
$summary.Activate() | Out-Null
$summary.Range("A1").Select() | Out-Null
